# Atualização automática 2025-09-30 15:32:05
# Adds the new SAP transaction entry "ZFI127" to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing table ends at row 72 (A72/B72 = "Consulta pagamentos por
# Fornecedor" / "FBL1N"). Append the new row right below it.
$ws.Range("A73").Value = "Extração Saldo de Razão das Contas Resumido"
$ws.Range("B73").Value = "ZFI127"

# Move the active selection past the newly added row, matching the
# post-edit cursor position recorded in the workbook.
[void]$ws.Range("A74").Select()
